$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose content actually differs row-to-row in this sheet.
# Everything else (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT,
# AW, AX, AY, ...) is identical for every data row and is left untouched.
$cols = @("A","B","E","F","G","H","K","L","M","N","Q","R","AC")

# The data rows (2..13) got reshuffled: each destination row's varying
# columns now hold what used to live in a different source row. Mapping is
# destination row -> source row (both are original, i.e. "before", row
# numbers).
$mapping = @{
    2  = 10
    3  = 2
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 11
    11 = 12
    12 = 13
    13 = 9
}

# Snapshot every source row's values *before* writing anything, since
# several destinations read from rows that are themselves about to be
# overwritten.
$old = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $old[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $old[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
